# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.796.49'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.734.55'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5181'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2747'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.32'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06123'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.736.70'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07051'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.89'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6373'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.521'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '76.67'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9996'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.803.87'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.43'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006627'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.957.17'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.170'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.723'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.146'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.46'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.507'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.05'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.766'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.20'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08273'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.684'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.498'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04464'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.610'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9739'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6131'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.668'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01569'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9996'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.905'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.80'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3816'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7218'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05355'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1125'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.170'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.02'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.88'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.581'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.38%  '
